$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating "2022-Q3" (so it
#    inherits the same layout/styles), placing it right before
#    "2022-Q3", then updating its fund figures.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Force these to be stored as text (matching the source data, which
# keeps these figures as strings rather than numbers).
$q4.Range("D2").Value = "'4.05"
$q4.Range("E2").Value = "'90.08"
$q4.Range("F2").Value = "'5.09"
$q4.Range("G2").Value = "'0.2061"

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new row 2 for the
#    2022-Q4 figures, shifting the existing 2022-Q3 / 2022-Q2 rows
#    down by one.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()

# Newly inserted row inherits the header row's formatting; strip it
# back to the default (unstyled) look used by the other data rows.
$total.Range("B2:D2").ClearFormats()

# Give A2 the same style as the other index cells (A3/A4).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.21

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.16

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q2"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.07000000000000001
